$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (17) into the two new
# rows, then strip the numeric-format style back off column F (the new rows'
# CropID column is left on the default/general format, unlike row 17's).
$ws.Range("B17:K17").Copy($ws.Range("B18:K18"))
$ws.Range("B17:K17").Copy($ws.Range("B19:K19"))
$ws.Cells.Item(18, 6).ClearFormats()
$ws.Cells.Item(19, 6).ClearFormats()

# Row 18: Item "나무" (wood)
$ws.Cells.Item(18, 2).Value = 20001001
$ws.Cells.Item(18, 3).Value = "나무"
$ws.Cells.Item(18, 4).Value = "나무를 베서 나온 나무조각. 건물을 짓거나 제작할 때 사용된다."
$ws.Cells.Item(18, 5).Value = "Item"
$ws.Cells.Item(18, 6).Value = -1
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 1000
$ws.Cells.Item(18, 10).Value = "Item/"
$ws.Cells.Item(18, 11).Value = "wood"

# Row 19: Item "돌" (stone)
$ws.Cells.Item(19, 2).Value = 20001002
$ws.Cells.Item(19, 3).Value = "돌"
$ws.Cells.Item(19, 4).Value = "돌을 캐서 나온 돌조각. 건물을 짓거나 제작할 때 사용된다."
$ws.Cells.Item(19, 5).Value = "Item"
$ws.Cells.Item(19, 6).Value = -1
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 1000
$ws.Cells.Item(19, 10).Value = "Item/"
$ws.Cells.Item(19, 11).Value = "stone"

# Update the saved selection to match the author's final view.
$ws.Range("K23").Select()
